$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2025-11-16 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-11-17 Monday", 2)

# New values for the 20x5 table of arithmetic problems, in row-major order
$newValues = @(
    "90-48=", "15+63=", "41+24=", "86+1=", "40-4=", "32+0=", "37-21=", "64-48=", "15-8=", "53-36=",
    "84+14=", "44-34=", "62-8=", "39-39=", "68-66=", "91-60=", "17+53=", "69-62=", "13+41=", "25-1=",
    "19-3=", "8+57=", "91-36=", "76-67=", "71-5=", "24-13=", "32-21=", "79-56=", "47+46=", "0+38=",
    "22+71=", "73+0=", "27+49=", "86-76=", "84-79=", "56-33=", "98-56=", "85-24=", "94-28=", "32+53=",
    "86-39=", "39+33=", "21+59=", "38-24=", "53+14=", "47-19=", "64+13=", "57+38=", "53-47=", "92-82=",
    "89-89=", "30+37=", "41+8=", "84-12=", "68+6=", "16+81=", "22+16=", "57+23=", "29+0=", "78-35=",
    "18+5=", "13+40=", "88-55=", "54-26=", "34-7=", "77-39=", "60-21=", "21-20=", "78-71=", "45-26=",
    "56-25=", "34+20=", "23+44=", "94-83=", "54-8=", "49-42=", "32+54=", "91-31=", "13+45=", "61+9=",
    "90-86=", "60+10=", "62-17=", "44+0=", "2+54=", "29-14=", "40+33=", "82-0=", "89-33=", "79-36=",
    "95-39=", "59-22=", "39+46=", "8+0=", "23+61=", "89-22=", "76-5=", "17+5=", "15+42=", "31+20="
)

$table = $d.Tables.Item(1)
$numRows = $table.Rows.Count
$numCols = $table.Columns.Count

$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $table.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated date and $idx table cells"
